# Apply updates described in the diff for InLrPt_00_Results_Summary.xlsx
$wb = $excel.ActiveWorkbook

# --- Sheet "Stats": update row 12 (C12 and G12) ---
$statsSheet = $wb.Worksheets.Item("Stats")
$statsSheet.Range("C12").Value = 0.714
$statsSheet.Range("G12").Value = 1

# --- Sheet "VIF": update VIF values in column C, rows 2-12 ---
$vifSheet = $wb.Worksheets.Item("VIF")
$vifSheet.Range("C2").Value  = 8.028118934152998
$vifSheet.Range("C3").Value  = 4.105928322727145
$vifSheet.Range("C4").Value  = 5.479306339060577
$vifSheet.Range("C5").Value  = 5.758439174509814
$vifSheet.Range("C6").Value  = 2.68855644444936
$vifSheet.Range("C7").Value  = 3.349120513095468
$vifSheet.Range("C8").Value  = 3.196186622121956
$vifSheet.Range("C9").Value  = 105.7667082703045
$vifSheet.Range("C10").Value = 135.4044845455055
$vifSheet.Range("C11").Value = 25.86720161629629
$vifSheet.Range("C12").Value = 3.942859723556538
